# Auto-generated edit script: update Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.272.67"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +6.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.115.59"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.47%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.15"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.16%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.107.01"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.59%  "

$ws.Range("E10").Value = "  +13.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +7.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +8.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.60"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.66%  "

$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.627.78"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.18"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.173.47"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.109.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.73%  "

$ws.Range("E21").Value = "  +4.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.727"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +7.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.11"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.78%  "

$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.56"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +10.79%  "

$ws.Range("E28").Value = "  +1.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.68"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +10.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.94"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.76%  "

$ws.Range("E33").Value = "  +3.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0871"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +13.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.41"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +16.19%  "

$ws.Range("E36").Value = "  +6.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.10"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.64%  "

$ws.Range("E38").Value = "  +20.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.68"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "439.40"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +9.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.73"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.916.52"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0370"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.279"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +11.51%  "

$ws.Range("E45").Value = "  +7.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +8.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.03"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.54%  "

$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.65"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.90%  "
